$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Mark Good Friday (row 7) and Washington's Birthday (row 13) as market-closed
# days, consistent with the other holiday rows that already carry a -1 in
# column B ("vol").
$ws.Range("B7").Value = -1
$ws.Range("B13").Value = -1

# Move the active selection to F11, matching the author's final cursor
# position when they saved the workbook.
$ws.Range("F11").Select()
